$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 61, shifting existing rows 61-143 down to 62-144.
$ws.Rows.Item(61).Insert()

# Populate the newly inserted row 61 with the new data record.
$ws.Cells.Item(61, 1).Value = 11
$ws.Cells.Item(61, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(61, 3).Value = "Bíobío"
$ws.Cells.Item(61, 4).Value = 44895
$ws.Cells.Item(61, 5).Value = 8
$ws.Cells.Item(61, 6).Value = "Fruta"
$ws.Cells.Item(61, 7).Value = 100108
$ws.Cells.Item(61, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(61, 9).Value = 100108002
$ws.Cells.Item(61, 10).Value = "Mango"
$ws.Cells.Item(61, 11).Value = "Sin especificar"
$ws.Cells.Item(61, 12).Value = "Primera"
$ws.Cells.Item(61, 13).Value = 100
$ws.Cells.Item(61, 14).Value = 7500
$ws.Cells.Item(61, 15).Value = 8000
$ws.Cells.Item(61, 16).Value = 7750
$ws.Cells.Item(61, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(61, 18).Value = "Perú"
$ws.Cells.Item(61, 19).Value = 1938
$ws.Cells.Item(61, 20).Value = 4
